# Updated symbol list on Thu Jan 26 14:53:16 UTC 2023 with GitHub Actions
# Refresh coin Price (D) and Volume(1h) (E) columns with latest scraped values.
# Values are plain text (not numbers) to preserve exact literal formatting
# (e.g. trailing zeros, scientific-looking small decimals, percent signs).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "307.40"
Set-TextValue $ws.Range("E2") "2.26%"

Set-TextValue $ws.Range("D3") "35.81"
Set-TextValue $ws.Range("E3") "1.26%"

Set-TextValue $ws.Range("D4") "5.082"
Set-TextValue $ws.Range("E4") "0.88%"

Set-TextValue $ws.Range("D5") "0.08065"
Set-TextValue $ws.Range("E5") "1.04%"

Set-TextValue $ws.Range("D6") "1.948"
Set-TextValue $ws.Range("E6") "1.03%"

Set-TextValue $ws.Range("D7") "4.161"
Set-TextValue $ws.Range("E7") "2.70%"

Set-TextValue $ws.Range("D8") "7.823"
Set-TextValue $ws.Range("E8") "0.44%"

Set-TextValue $ws.Range("D9") "0.9342"
Set-TextValue $ws.Range("E9") "1.16%"

Set-TextValue $ws.Range("D10") "0.1336"
Set-TextValue $ws.Range("E10") "-10.49%"

Set-TextValue $ws.Range("D11") "0.1913"
Set-TextValue $ws.Range("E11") "1.03%"

Set-TextValue $ws.Range("D12") "0.09243"
Set-TextValue $ws.Range("E12") "0.55%"

Set-TextValue $ws.Range("D13") "0.03522"
Set-TextValue $ws.Range("E13") "2.54%"

Set-TextValue $ws.Range("D14") "0.09890"
Set-TextValue $ws.Range("E14") "0.11%"

Set-TextValue $ws.Range("D15") "0.001433"
Set-TextValue $ws.Range("E15") "3.08%"

Set-TextValue $ws.Range("D16") "0.005868"
Set-TextValue $ws.Range("E16") "1.85%"

Set-TextValue $ws.Range("E17") "2.68%"

Set-TextValue $ws.Range("D18") "2.920"
Set-TextValue $ws.Range("E18") "-2.37%"

Set-TextValue $ws.Range("D19") "0.3459"
Set-TextValue $ws.Range("E19") "1.63%"

Set-TextValue $ws.Range("D20") "0.1336"
Set-TextValue $ws.Range("E20") "3.30%"

Set-TextValue $ws.Range("D21") "5.170"
Set-TextValue $ws.Range("E21") "2.42%"

Set-TextValue $ws.Range("D22") "0.2628"
Set-TextValue $ws.Range("E22") "9.34%"

Set-TextValue $ws.Range("D23") "0.04401"
Set-TextValue $ws.Range("E23") "-1.36%"

Set-TextValue $ws.Range("D24") "0.001241"
Set-TextValue $ws.Range("E24") "2.01%"

Set-TextValue $ws.Range("D25") "0.004762"
Set-TextValue $ws.Range("E25") "-0.27%"

Set-TextValue $ws.Range("D26") "0.0001302"
Set-TextValue $ws.Range("E26") "5.59%"

Set-TextValue $ws.Range("D27") "0.0003141"
Set-TextValue $ws.Range("E27") "4.42%"

Set-TextValue $ws.Range("D39") "0.01984"
Set-TextValue $ws.Range("E39") "4.12%"

Set-TextValue $ws.Range("D40") "0.05013"
Set-TextValue $ws.Range("E40") "6.16%"

Set-TextValue $ws.Range("D41") "0.01122"

Set-TextValue $ws.Range("D42") "0.007622"
Set-TextValue $ws.Range("E42") "3.32%"

Set-TextValue $ws.Range("D43") "0.1371"
Set-TextValue $ws.Range("E43") "3.19%"

Set-TextValue $ws.Range("D44") "0.002103"
Set-TextValue $ws.Range("E44") "-0.58%"

Set-TextValue $ws.Range("D45") "0.01135"
Set-TextValue $ws.Range("E45") "21.66%"

Set-TextValue $ws.Range("D46") "0.00006404"
Set-TextValue $ws.Range("E46") "2.26%"

Set-TextValue $ws.Range("D47") "0.00000000751"
Set-TextValue $ws.Range("E47") "-0.07%"

Set-TextValue $ws.Range("D48") "64.96"
Set-TextValue $ws.Range("E48") "-0.07%"

Set-TextValue $ws.Range("D49") "0.001195"
Set-TextValue $ws.Range("E49") "-28.14%"

Set-TextValue $ws.Range("D50") "0.00002104"
Set-TextValue $ws.Range("E50") "-0.07%"

Set-TextValue $ws.Range("D51") "0.0002004"
Set-TextValue $ws.Range("E51") "-0.07%"
